# Auto-generated edit script applying the cryptos.xlsx data refresh diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.418.16'
$ws.Range('E2').Value = '  -1.70%  '
$ws.Range('D3').Value = '3.483.38'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.57'
$ws.Range('E5').Value = '  -2.40%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '173.54'
$ws.Range('E6').Value = '  -3.38%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '3.486.83'
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.594'
$ws.Range('E9').Value = '  -3.32%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.132'
$ws.Range('E10').Value = '  -4.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.89'
$ws.Range('E11').Value = '  -1.88%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.420'
$ws.Range('E12').Value = '  -3.55%  '
$ws.Range('D13').Value = '4.108.21'
$ws.Range('E13').Value = '  +0.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '30.13'
$ws.Range('E14').Value = '  -5.85%  '
$ws.Range('E15').Value = '  -1.49%  '
$ws.Range('D16').Value = '66.549.67'
$ws.Range('E16').Value = '  -1.54%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000171'
$ws.Range('E17').Value = '  -3.85%  '
$ws.Range('D18').Value = '3.498.99'
$ws.Range('E18').Value = '  +0.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.99'
$ws.Range('E19').Value = '  -5.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.77'
$ws.Range('E20').Value = '  -3.49%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '379.17'
$ws.Range('E21').Value = '  -2.97%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.80'
$ws.Range('E22').Value = '  -2.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.545'
$ws.Range('E23').Value = '  +0.65%  '
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '72.03'
$ws.Range('E26').Value = '  -1.43%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000118'
$ws.Range('E27').Value = '  -3.47%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.79'
$ws.Range('E28').Value = '  -5.31%  '
$ws.Range('E29').Value = '  -1.02%  '
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '24.34'
$ws.Range('E31').Value = '  +3.24%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.84'
$ws.Range('E32').Value = '  -5.11%  '
$ws.Range('E33').Value = '  -3.44%  '
$ws.Range('B34').Value = 'USDe'
$ws.Range('C34').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('B35').Value = 'Fetch.AI'
$ws.Range('C35').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.31'
$ws.Range('E35').Value = '  -7.63%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.22'
$ws.Range('E36').Value = '  -2.54%  '
$ws.Range('E37').Value = '  -2.89%  '
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '159.85'
$ws.Range('E38').Value = '  -1.84%  '
$ws.Range('B39').Value = 'Mantle'
$ws.Range('C39').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.891'
$ws.Range('E39').Value = '  +1.09%  '
$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '29.07'
$ws.Range('E40').Value = '  +10.81%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.77'
$ws.Range('E41').Value = '  -6.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.51'
$ws.Range('E42').Value = '  -2.22%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.40'
$ws.Range('E43').Value = '  -6.09%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.52'
$ws.Range('E44').Value = '  -10.62%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '2.672.01'
$ws.Range('E45').Value = '  -5.56%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0692'
$ws.Range('E46').Value = '  -4.65%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '40.59'
$ws.Range('E47').Value = '  -2.38%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '24.32'
$ws.Range('E48').Value = '  -9.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0290'
$ws.Range('E49').Value = '  -3.49%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '313.38'
$ws.Range('E50').Value = '  -6.33%  '
$ws.Range('E51').Value = '  -4.61%  '
